# Insert a new weekly price record at row 96, pushing the existing
# rows 96-186 down to 97-187 (dimension grows from R186 to R187).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new observation.
$ws.Range("A96").Value = 10
$ws.Range("B96").Value = "Vega Modelo de Temuco"
$ws.Range("C96").Value = "La Araucanía"
$ws.Range("D96").Value = 44539
$ws.Range("E96").Value = 9
$ws.Range("F96").Value = 100112039
$ws.Range("G96").Value = "Ciboulette"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 65
$ws.Range("K96").Value = 5000
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = 5000
$ws.Range("N96").Value = "$/docena de atados"
$ws.Range("O96").Value = "Provincia de Cautín"
$ws.Range("P96").Value = 1667
$ws.Range("Q96").Value = 3
$ws.Range("R96").Value = "Hortaliza"
